$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the date-ordered data block (row 542),
# shifting the existing rows 542:588 down to 545:591.
$ws.Rows("542:544").Insert()

# Common fields shared by every data row in this sheet.
$colA = 8
$colB = "Terminal La Palmera de La Serena"
$colC = "Coquimbo"
$colE = 4
$colF = "Fruta"
$colG = 100101
$colH = "Berries"
$colI = 100101007
$colJ = "Kiwi"
$colK = "Hayward"
$colQ = "$/bins (450 kilos)"
$colR = "Región de O'Higgins"
$newDate = 45106
$colT = 450

# Row 542: Especial
$r = 542
$ws.Cells.Item($r, 1).Value2 = $colA
$ws.Cells.Item($r, 2).Value2 = $colB
$ws.Cells.Item($r, 3).Value2 = $colC
$ws.Cells.Item($r, 4).Value2 = $newDate
$ws.Cells.Item($r, 5).Value2 = $colE
$ws.Cells.Item($r, 6).Value2 = $colF
$ws.Cells.Item($r, 7).Value2 = $colG
$ws.Cells.Item($r, 8).Value2 = $colH
$ws.Cells.Item($r, 9).Value2 = $colI
$ws.Cells.Item($r, 10).Value2 = $colJ
$ws.Cells.Item($r, 11).Value2 = $colK
$ws.Cells.Item($r, 12).Value2 = "Especial"
$ws.Cells.Item($r, 13).Value2 = 10
$ws.Cells.Item($r, 14).Value2 = 350000
$ws.Cells.Item($r, 15).Value2 = 360000
$ws.Cells.Item($r, 16).Value2 = 355000
$ws.Cells.Item($r, 17).Value2 = $colQ
$ws.Cells.Item($r, 18).Value2 = $colR
$ws.Cells.Item($r, 19).Value2 = 789
$ws.Cells.Item($r, 20).Value2 = $colT

# Row 543: Primera
$r = 543
$ws.Cells.Item($r, 1).Value2 = $colA
$ws.Cells.Item($r, 2).Value2 = $colB
$ws.Cells.Item($r, 3).Value2 = $colC
$ws.Cells.Item($r, 4).Value2 = $newDate
$ws.Cells.Item($r, 5).Value2 = $colE
$ws.Cells.Item($r, 6).Value2 = $colF
$ws.Cells.Item($r, 7).Value2 = $colG
$ws.Cells.Item($r, 8).Value2 = $colH
$ws.Cells.Item($r, 9).Value2 = $colI
$ws.Cells.Item($r, 10).Value2 = $colJ
$ws.Cells.Item($r, 11).Value2 = $colK
$ws.Cells.Item($r, 12).Value2 = "Primera"
$ws.Cells.Item($r, 13).Value2 = 10
$ws.Cells.Item($r, 14).Value2 = 310000
$ws.Cells.Item($r, 15).Value2 = 320000
$ws.Cells.Item($r, 16).Value2 = 315000
$ws.Cells.Item($r, 17).Value2 = $colQ
$ws.Cells.Item($r, 18).Value2 = $colR
$ws.Cells.Item($r, 19).Value2 = 700
$ws.Cells.Item($r, 20).Value2 = $colT

# Row 544: Segunda
$r = 544
$ws.Cells.Item($r, 1).Value2 = $colA
$ws.Cells.Item($r, 2).Value2 = $colB
$ws.Cells.Item($r, 3).Value2 = $colC
$ws.Cells.Item($r, 4).Value2 = $newDate
$ws.Cells.Item($r, 5).Value2 = $colE
$ws.Cells.Item($r, 6).Value2 = $colF
$ws.Cells.Item($r, 7).Value2 = $colG
$ws.Cells.Item($r, 8).Value2 = $colH
$ws.Cells.Item($r, 9).Value2 = $colI
$ws.Cells.Item($r, 10).Value2 = $colJ
$ws.Cells.Item($r, 11).Value2 = $colK
$ws.Cells.Item($r, 12).Value2 = "Segunda"
$ws.Cells.Item($r, 13).Value2 = 20
$ws.Cells.Item($r, 14).Value2 = 260000
$ws.Cells.Item($r, 15).Value2 = 270000
$ws.Cells.Item($r, 16).Value2 = 265000
$ws.Cells.Item($r, 17).Value2 = $colQ
$ws.Cells.Item($r, 18).Value2 = $colR
$ws.Cells.Item($r, 19).Value2 = 589
$ws.Cells.Item($r, 20).Value2 = $colT
